# Contest 22 DC vs RCB.
# Fill in the results row for Contest 22 (row 34 on Sheet1) with each
# participant's score for this match. Column C34 already holds the
# "DC vs RCB" label (shared string); the per-player input cells
# (E, H, K, N, Q, T, W, Z, AC) were still blank and are populated here.
# The adjoining "format"/rank formulas in D, G, J, M, P, S, V, Y, AB and
# the Total row (45) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E34").Value  = 70    # Jaya
$ws.Range("H34").Value  = 100   # Justin
$ws.Range("K34").Value  = 20    # Ram
$ws.Range("N34").Value  = 60    # Sibi
$ws.Range("Q34").Value  = 30    # Sundar
$ws.Range("T34").Value  = 80    # Balaji
$ws.Range("W34").Value  = 50    # Upili
$ws.Range("Z34").Value  = 40    # Vicky
$ws.Range("AC34").Value = 0     # Raghu
